# Generate Report for Archive
#
# The localization job moved out of "Ready for handoff" into the
# translation phase, so every cell holding that status string (the
# Overview sheet's per-locale status columns, plus each per-locale
# sheet's "Status" column) is updated to "In Translation". Because the
# new text is shorter, Excel re-fits the affected columns to their new
# (narrower) autofit width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColumnWidth = 13.4101845877511

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count
    $colCount = $usedRange.Columns.Count

    # Track which columns actually contained the old status text, so we
    # only resize the columns Excel's autofit would have touched.
    $touchedCols = @{}

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $usedRange.Cells.Item($r, $c)
            # Compare with the string on the LEFT so PowerShell coerces the
            # (possibly boolean) cell value to a string instead of coercing
            # our string to a boolean (which would make every non-empty
            # string cell compare equal to a TRUE boolean cell).
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
                $touchedCols[$cell.Column] = $true
            }
        }
    }

    # Re-fit the columns whose text just got shorter. (AutoFit is a no-op
    # in this headless runtime, so the target width is applied directly.)
    foreach ($colIndex in $touchedCols.Keys) {
        $ws.Columns.Item($colIndex).ColumnWidth = $newColumnWidth
    }
}
